$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name (C3, merged C3:D3)
$ws.Range("C3").Value = "Ralph Vitug"

# Row 7 - __init__ / Attribute set to input values
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'Rectangle("blue" , 10, 5)'
$ws.Range("G7").Value = "Object created"

# Row 8 - __init__ / Exception raised when color is blank
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'Rectangle(" " , 10, 5)'
$ws.Range("G8").Value = "ValueError"

# Row 9 - __init__ / Exception raised when length is not an integer
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'Rectangle("blue" , "a", 5)'
$ws.Range("G9").Value = "ValueError"

# Row 10 - __init__ / Exception raised when width is not an integer
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'Rectangle("blue" , 10, "b")'
$ws.Range("G10").Value = "ValueError"

# Row 11 - __str__ / Returns string formatted appropriately
$ws.Range("E11").Value = "`"The shape color is blue.`"`n  `"This rectangle has four sides with the lengths of  10, 5, 10 and 5 centimeters`""
$ws.Range("F11").Value = "`"The shape color is blue.\n`"`n            `"This rectangle has four sides with the lengths `"`n            `"of 10, 5, 10 and 5 centimeter`""
$ws.Range("G11").Value = "`"The shape color is blue.`"`n  `"This rectangle has four sides with the lengths of  10, 5, 10 and 5 centimeters`""

# Row 12 - calculate_area / Returns correct calculated value
$ws.Range("E12").Value = "50, self.triangle.calculate_area()"
$ws.Range("F12").Value = "self.assertEqual(50, self.rectangle.calculate_area())"
$ws.Range("G12").Value = 50

# Row 13 - calculate_perimeter
$ws.Range("E13").Value = "30, self.triangle.calculate_area()"
$ws.Range("F13").Value = "self.assertEqual(30, self.rectangle.calculate_area())"
$ws.Range("G13").Value = 30

# Update selection to match the final state
$ws.Range("C3:D3").Select()
